# ---------------------------------------------------------------------------
# Update with Correct Forecast output
#
# 1) Rename Sheet1 -> "Sales vs PO" and insert a new "Order Week" column
#    (shifting the old PO_Requested_Qty column from C to D), refreshing the
#    ds/y/Order Week/PO_Requested_Qty values to the corrected forecast.
# 2) Add "Weekly Growth" sheet with the weekly PO growth figures.
# 3) Add "Volume Insights" sheet with aggregate PO stats.
# 4) Add "Prediction Info" sheet with the next-week PO prediction.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Sales vs PO" -------------------------------------------------
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sales vs PO"

# Insert a new column before the existing PO_Requested_Qty column (C),
# pushing it to column D while preserving its header/value formatting.
$ws1.Columns.Item(3).Insert()

$ws1.Range("C1").Value = "Order Week"

$salesPoRows = @(
    # A (ds)   B (y)  C (Order Week)  D (PO_Requested_Qty)
    @(45571,   0,     45565,          0),
    @(45578,   0,     45572,          0),
    @(45585,   0,     45579,          0),
    @(45592,   0,     45586,          0),
    @(45599,   0,     45593,          0),
    @(45606,   0,     45600,          0),
    @(45613,   0,     45607,          0),
    @(45620,   11,    45614,          0),
    @(45627,   6,     45621,          0),
    @(45634,   9,     45628,          0),
    @(45641,   3,     45635,          0),
    @(45648,   8,     45642,          0),
    @(45655,   6,     45649,          0)
)

$r = 2
foreach ($row in $salesPoRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 3).NumberFormat = $ws1.Cells.Item($r, 1).NumberFormat
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Header style (bold, bordered, centered) used throughout the workbook.
$headerCell = $ws1.Range("A1")

function Set-PageMargins($sheet) {
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# --- Sheet 2: "Weekly Growth" ----------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"
Set-PageMargins $ws2

$headerCell.Copy($ws2.Range("A1"))
$ws2.Range("A1").Value = "ds"
$headerCell.Copy($ws2.Range("B1"))
$ws2.Range("B1").Value = "PO_Requested_Qty"
$headerCell.Copy($ws2.Range("C1"))
$ws2.Range("C1").Value = "Growth%"

$weeklyGrowthRows = @(
    @(45572, 20,  0),
    @(45586, 40,  100),
    @(45607, 180, 350)
)

$r = 2
foreach ($row in $weeklyGrowthRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- Sheet 3: "Volume Insights" ---------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"
Set-PageMargins $ws3

$headerCell.Copy($ws3.Range("A1"))
$ws3.Range("A1").Value = "Total_PO_Quantity"
$headerCell.Copy($ws3.Range("B1"))
$ws3.Range("B1").Value = "Average_PO_Quantity"
$headerCell.Copy($ws3.Range("C1"))
$ws3.Range("C1").Value = "Max_PO_Quantity"
$headerCell.Copy($ws3.Range("D1"))
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 240
$ws3.Range("B2").Value = 80
$ws3.Range("C2").Value = 180
$ws3.Range("D2").Value = 20

# --- Sheet 4: "Prediction Info" ---------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"
Set-PageMargins $ws4

$headerCell.Copy($ws4.Range("A1"))
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"

$ws4.Range("A2").Value = 239.9999999999999

# --- Activate first sheet to match original view ---------------------------
$ws1.Activate()
